$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: forecast week starting dates (serial date numbers)
$ws.Range("B1").Value = 45627
$ws.Range("C1").Value = 45634
$ws.Range("D1").Value = 45641
$ws.Range("E1").Value = 45648
$ws.Range("F1").Value = 45655
$ws.Range("G1").Value = 45662
$ws.Range("H1").Value = 45669
$ws.Range("I1").Value = 45676
$ws.Range("J1").Value = 45683
$ws.Range("K1").Value = 45690
$ws.Range("L1").Value = 45697

# Rows 2-12: VCI3M forecast values
# Row 2
$ws.Range("B2").Value = 63.7
$ws.Range("C2").Value = 60.5
$ws.Range("D2").Value = 57.8
$ws.Range("E2").Value = 55.9
$ws.Range("F2").Value = 54.8
$ws.Range("G2").Value = 54.5
$ws.Range("H2").Value = 54.9
$ws.Range("I2").Value = 56
$ws.Range("J2").Value = 57.4
$ws.Range("K2").Value = 58.9
$ws.Range("L2").Value = 60.4

# Row 3
$ws.Range("B3").Value = 63.2
$ws.Range("C3").Value = 59.4
$ws.Range("D3").Value = 56.3
$ws.Range("E3").Value = 54.2
$ws.Range("F3").Value = 53
$ws.Range("G3").Value = 52.8
$ws.Range("H3").Value = 53.3
$ws.Range("I3").Value = 54.4
$ws.Range("J3").Value = 55.9
$ws.Range("K3").Value = 57.5
$ws.Range("L3").Value = 59.1

# Row 4
$ws.Range("B4").Value = 63.7
$ws.Range("C4").Value = 61.9
$ws.Range("D4").Value = 59.9
$ws.Range("E4").Value = 57.8
$ws.Range("F4").Value = 55.7
$ws.Range("G4").Value = 53.7
$ws.Range("H4").Value = 51.6
$ws.Range("I4").Value = 49.6
$ws.Range("J4").Value = 47.6
$ws.Range("K4").Value = 45.6
$ws.Range("L4").Value = 43.8

# Row 5
$ws.Range("B5").Value = 66.59999999999999
$ws.Range("C5").Value = 65.2
$ws.Range("D5").Value = 64.8
$ws.Range("E5").Value = 65.2
$ws.Range("F5").Value = 66.5
$ws.Range("G5").Value = 68.3
$ws.Range("H5").Value = 70.59999999999999
$ws.Range("I5").Value = 72.90000000000001
$ws.Range("J5").Value = 75.09999999999999
$ws.Range("K5").Value = 76.8
$ws.Range("L5").Value = 78.09999999999999

# Row 6
$ws.Range("B6").Value = 70.09999999999999
$ws.Range("C6").Value = 64.5
$ws.Range("D6").Value = 59.8
$ws.Range("E6").Value = 56.2
$ws.Range("F6").Value = 54.1
$ws.Range("G6").Value = 53.3
$ws.Range("H6").Value = 53.8
$ws.Range("I6").Value = 55.5
$ws.Range("J6").Value = 58.1
$ws.Range("K6").Value = 61.2
$ws.Range("L6").Value = 64.5

# Row 7
$ws.Range("B7").Value = 66.5
$ws.Range("C7").Value = 65.3
$ws.Range("D7").Value = 64.09999999999999
$ws.Range("E7").Value = 63
$ws.Range("F7").Value = 62
$ws.Range("G7").Value = 61.2
$ws.Range("H7").Value = 60.4
$ws.Range("I7").Value = 59.7
$ws.Range("J7").Value = 58.9
$ws.Range("K7").Value = 58
$ws.Range("L7").Value = 56.8

# Row 8
$ws.Range("B8").Value = 63.7
$ws.Range("C8").Value = 61.3
$ws.Range("D8").Value = 58.7
$ws.Range("E8").Value = 56
$ws.Range("F8").Value = 53.4
$ws.Range("G8").Value = 50.7
$ws.Range("H8").Value = 48.2
$ws.Range("I8").Value = 45.7
$ws.Range("J8").Value = 43.4
$ws.Range("K8").Value = 41.1
$ws.Range("L8").Value = 38.9

# Row 9
$ws.Range("B9").Value = 51.6
$ws.Range("C9").Value = 48.3
$ws.Range("D9").Value = 46.1
$ws.Range("E9").Value = 45.4
$ws.Range("F9").Value = 46.2
$ws.Range("G9").Value = 48.6
$ws.Range("H9").Value = 52.5
$ws.Range("I9").Value = 57.4
$ws.Range("J9").Value = 63
$ws.Range("K9").Value = 68.7
$ws.Range("L9").Value = 74.09999999999999

# Row 10
$ws.Range("B10").Value = 62.3
$ws.Range("C10").Value = 58.4
$ws.Range("D10").Value = 55.4
$ws.Range("E10").Value = 53.5
$ws.Range("F10").Value = 52.7
$ws.Range("G10").Value = 53.1
$ws.Range("H10").Value = 54.6
$ws.Range("I10").Value = 56.9
$ws.Range("J10").Value = 59.8
$ws.Range("K10").Value = 63
$ws.Range("L10").Value = 66.2

# Row 11
$ws.Range("B11").Value = 70.2
$ws.Range("C11").Value = 62.5
$ws.Range("D11").Value = 56
$ws.Range("E11").Value = 51.4
$ws.Range("F11").Value = 48.9
$ws.Range("G11").Value = 48.7
$ws.Range("H11").Value = 50.7
$ws.Range("I11").Value = 54.7
$ws.Range("J11").Value = 60.2
$ws.Range("K11").Value = 66.59999999999999
$ws.Range("L11").Value = 73.09999999999999

# Row 12
$ws.Range("B12").Value = 59.3
$ws.Range("C12").Value = 58.1
$ws.Range("D12").Value = 56.9
$ws.Range("E12").Value = 56
$ws.Range("F12").Value = 55.1
$ws.Range("G12").Value = 54.3
$ws.Range("H12").Value = 53.5
$ws.Range("I12").Value = 52.7
$ws.Range("J12").Value = 51.8
$ws.Range("K12").Value = 50.8
$ws.Range("L12").Value = 49.9
